# "save increase/decrease calc in excel"
# Add two new helper columns (K, L) on Sheet1 that compute the
# increase/decrease ratios already present as raw numbers in columns D-H:
#   K = E / D   (ref expr vs new expr change ratio)
#   L = H / F   (ref flux vs new flux change ratio)
# for every data row (2 through 98), then leave the selection on M5
# (the cell just to the right of the new columns), matching the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 98; $r++) {
    $ws.Range("K$r").Formula = "=E$r/D$r"
    $ws.Range("L$r").Formula = "=H$r/F$r"
}

[void]$ws.Range("M5").Select()
